# Overwrite the "current test run" row (row 2) of the TestData sheet with a
# freshly generated user's details - FirstName, LastName, Email, Password -
# mirroring the rotating scratch-row pattern already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Andrew"
$ws.Range("C2").Value = "Williams"
$ws.Range("D2").Value = "andrew.williams1735829998241@test.com"
$ws.Range("E2").Value = '5xX$H^{-t'
